# CPI corrected, new Plots
# Update the CSCC values in column C (rows 2-7) to reflect the corrected CPI figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 180.789157331046
$ws.Range("C3").Value = 420.5017506766929
$ws.Range("C4").Value = -19.94617240771698
$ws.Range("C5").Value = 71.43296813034826
$ws.Range("C6").Value = 22.14866148088394
$ws.Range("C7").Value = 4.527065827029297
